# Update the hardcoded "datetimeFigureOut" date field text from 28.08.2015
# to 31.08.2015 across the slide master and every slide layout's date
# placeholder ("Datumsplatzhalter").

$p = $ppt.ActivePresentation

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Datumsplatzhalter*") {
            if ($shp.TextFrame.TextRange.Text -eq "28.08.2015") {
                $shp.TextFrame.TextRange.Text = "31.08.2015"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShape $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShape $layout.Shapes
}
